# Daily attendance processing - 2025-12-30 21:31:19
# Reorders the "Recorded By" (column G) value on each row so that the
# literal token "System" (exact case) is moved from its current position
# to the end of the comma-separated list, while leaving every other
# token (including the unrelated lowercase "system") in its original
# relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $parts = $raw -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $systemIndex = $i
        }
    }

    if ($systemIndex -eq -1) {
        continue
    }

    if ($systemIndex -eq ($parts.Count - 1)) {
        continue
    }

    $newParts = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $systemIndex) {
            $newParts += $parts[$i]
        }
    }
    $newParts += "System"

    $cell.Value = $newParts -join ", "
}
